# Week 16 stat logging + season totals update (Packers Players Data.xlsx)
$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$rushing = $wb.Worksheets.Item("Rushing")

# Row 3: A.Jones
$rushing.Range("C3").Value = 107
$rushing.Range("D3").Value = 53
$rushing.Range("E3").Value = 8

# Row 4: A.Dillon
$rushing.Range("C4").Value = 80
$rushing.Range("D4").Value = 56
$rushing.Range("F4").Value = 32

# Row 7: A.Lazard
$rushing.Range("C7").Value = 2

# Row 9: E.St. Brown
$rushing.Range("C9").Value = 3

# --- Receiving sheet ---
$receiving = $wb.Worksheets.Item("Receiving")

# Row 2: A.Jones
$receiving.Range("C2").Value = 58
$receiving.Range("D2").Value = 48
$receiving.Range("G2").Value = 13
$receiving.Range("H2").Value = 10

# Row 3: A.Dillon
$receiving.Range("C3").Value = 26
$receiving.Range("D3").Value = 22
$receiving.Range("G3").Value = 4
$receiving.Range("H3").Value = 4

# Row 5: D.Adams
$receiving.Range("C5").Value = 120
$receiving.Range("D5").Value = 100
$receiving.Range("E5").Value = 39
$receiving.Range("F5").Value = 29
$receiving.Range("G5").Value = 25
$receiving.Range("H5").Value = 22

# Row 7: A.Lazard
$receiving.Range("C7").Value = 45
$receiving.Range("D7").Value = 34
$receiving.Range("E7").Value = 11
$receiving.Range("F7").Value = 7
$receiving.Range("G7").Value = 13
$receiving.Range("H7").Value = 6

# Row 11: J.Winfree
$receiving.Range("E11").Value = 1

# Row 13: M.Lewis
$receiving.Range("C13").Value = 26
$receiving.Range("D13").Value = 23
$receiving.Range("G13").Value = 5

# Row 14: J.Deguara
$receiving.Range("C14").Value = 22
$receiving.Range("D14").Value = 17
$receiving.Range("G14").Value = 4
$receiving.Range("H14").Value = 2
